$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 6848  # was 6822
$ws.Cells.Item(3, 10).Value = 7234  # was 7211
$ws.Cells.Item(4, 10).Value = 1576  # was 1572
$ws.Cells.Item(5, 10).Value = 569  # was 566
$ws.Cells.Item(6, 10).Value = 9660  # was 9621
$ws.Cells.Item(7, 10).Value = 25887  # was 25792
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 10).Value = 204  # was 203
$ws.Cells.Item(6, 10).Value = 196  # was 195
$ws.Cells.Item(7, 10).Value = 747  # was 746
$ws.Cells.Item(8, 10).Value = 1625  # was 1621
$ws.Cells.Item(10, 10).Value = 190  # was 189
$ws.Cells.Item(11, 10).Value = 447  # was 444
$ws.Cells.Item(13, 10).Value = 31  # was 30
$ws.Cells.Item(15, 10).Value = 308  # was 306
$ws.Cells.Item(19, 10).Value = 756  # was 753
$ws.Cells.Item(23, 10).Value = 240  # was 238
$ws.Cells.Item(24, 10).Value = 83  # was 81
$ws.Cells.Item(27, 10).Value = 153  # was 152
$ws.Cells.Item(29, 10).Value = 1404  # was 1398
$ws.Cells.Item(31, 10).Value = 261  # was 259
$ws.Cells.Item(33, 10).Value = 1169  # was 1164
$ws.Cells.Item(36, 10).Value = 353  # was 351
$ws.Cells.Item(37, 10).Value = 798  # was 796
$ws.Cells.Item(41, 10).Value = 180  # was 178
$ws.Cells.Item(42, 10).Value = 1115  # was 1110
$ws.Cells.Item(46, 10).Value = 88  # was 86
$ws.Cells.Item(47, 10).Value = 192  # was 191
$ws.Cells.Item(48, 10).Value = 295  # was 293
$ws.Cells.Item(49, 10).Value = 163  # was 162
$ws.Cells.Item(51, 10).Value = 316  # was 313
$ws.Cells.Item(52, 10).Value = 660  # was 654
$ws.Cells.Item(53, 10).Value = 379  # was 372
$ws.Cells.Item(55, 10).Value = 396  # was 394
$ws.Cells.Item(58, 10).Value = 17  # was 16
$ws.Cells.Item(60, 10).Value = 150  # was 149
$ws.Cells.Item(63, 10).Value = 79  # was 80
$ws.Cells.Item(64, 10).Value = 172  # was 170
$ws.Cells.Item(67, 10).Value = 966  # was 967
$ws.Cells.Item(72, 10).Value = 100  # was 99
$ws.Cells.Item(73, 10).Value = 250  # was 248
$ws.Cells.Item(75, 10).Value = 79  # was 78
$ws.Cells.Item(78, 10).Value = 307  # was 305
$ws.Cells.Item(79, 10).Value = 730  # was 728
$ws.Cells.Item(83, 10).Value = 516  # was 514
$ws.Cells.Item(85, 10).Value = 1064  # was 1063
$ws.Cells.Item(90, 10).Value = 276  # was 274
$ws.Cells.Item(93, 10).Value = 107  # was 106
$ws.Cells.Item(95, 10).Value = 378  # was 376
$ws.Cells.Item(96, 10).Value = 281  # was 280
$ws.Cells.Item(97, 10).Value = 238  # was 237
$ws.Cells.Item(98, 10).Value = 192  # was 188
$ws.Cells.Item(99, 10).Value = 397  # was 396
$ws.Cells.Item(101, 10).Value = 25887  # was 25792
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(4, 10).Value = 18  # was 17
$ws.Cells.Item(7, 10).Value = 281  # was 280
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 232  # was 231
$ws.Cells.Item(7, 10).Value = 747  # was 746
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(4, 10).Value = 26  # was 25
$ws.Cells.Item(6, 10).Value = 204  # was 202
$ws.Cells.Item(7, 10).Value = 447  # was 444
$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 10).Value = 373  # was 372
$ws.Cells.Item(7, 10).Value = 1064  # was 1063
$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 10).Value = 155  # was 154
$ws.Cells.Item(3, 10).Value = 189  # was 185
$ws.Cells.Item(6, 10).Value = 281  # was 280
$ws.Cells.Item(7, 10).Value = 660  # was 654
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(4, 10).Value = 15  # was 14
$ws.Cells.Item(5, 10).Value = 2  # was 1
$ws.Cells.Item(6, 10).Value = 251  # was 246
$ws.Cells.Item(7, 10).Value = 379  # was 372
$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 10).Value = 437  # was 436
$ws.Cells.Item(3, 10).Value = 485  # was 484
$ws.Cells.Item(6, 10).Value = 574  # was 572
$ws.Cells.Item(7, 10).Value = 1625  # was 1621
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 10).Value = 193  # was 192
$ws.Cells.Item(6, 10).Value = 137  # was 136
$ws.Cells.Item(7, 10).Value = 516  # was 514
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 10).Value = 271  # was 269
$ws.Cells.Item(3, 10).Value = 389  # was 388
$ws.Cells.Item(6, 10).Value = 407  # was 405
$ws.Cells.Item(7, 10).Value = 1169  # was 1164
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 10).Value = 131  # was 129
$ws.Cells.Item(7, 10).Value = 378  # was 376
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 10).Value = 236  # was 235
$ws.Cells.Item(3, 10).Value = 268  # was 267
$ws.Cells.Item(7, 10).Value = 798  # was 796
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 10).Value = 108  # was 107
$ws.Cells.Item(7, 10).Value = 397  # was 396
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 10).Value = 92  # was 91
$ws.Cells.Item(6, 10).Value = 86  # was 85
$ws.Cells.Item(7, 10).Value = 261  # was 259
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 10).Value = 359  # was 358
$ws.Cells.Item(6, 10).Value = 267  # was 269
$ws.Cells.Item(7, 10).Value = 966  # was 967
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(4, 10).Value = 12  # was 11
$ws.Cells.Item(7, 10).Value = 163  # was 162
$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 10).Value = 425  # was 424
$ws.Cells.Item(3, 10).Value = 498  # was 496
$ws.Cells.Item(4, 10).Value = 73  # was 72
$ws.Cells.Item(5, 10).Value = 53  # was 52
$ws.Cells.Item(6, 10).Value = 355  # was 354
$ws.Cells.Item(7, 10).Value = 1404  # was 1398
$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(3, 10).Value = 55  # was 54
$ws.Cells.Item(6, 10).Value = 144  # was 143
$ws.Cells.Item(7, 10).Value = 295  # was 293
$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 10).Value = 182  # was 181
$ws.Cells.Item(5, 10).Value = 29  # was 27
$ws.Cells.Item(7, 10).Value = 756  # was 753
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(6, 10).Value = 71  # was 70
$ws.Cells.Item(7, 10).Value = 196  # was 195
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(6, 10).Value = 105  # was 103
$ws.Cells.Item(7, 10).Value = 180  # was 178
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 10).Value = 239  # was 237
$ws.Cells.Item(3, 10).Value = 221  # was 220
$ws.Cells.Item(6, 10).Value = 588  # was 586
$ws.Cells.Item(7, 10).Value = 1115  # was 1110
$ws = $wb.Worksheets.Item('Boystown')
$ws.Cells.Item(5, 10).Value = 15  # was 14
$ws.Cells.Item(6, 10).Value = 31  # was 30
$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(2, 10).Value = 43  # was 42
$ws.Cells.Item(7, 10).Value = 190  # was 189
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(3, 10).Value = 96  # was 94
$ws.Cells.Item(7, 10).Value = 307  # was 305
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 10).Value = 78  # was 77
$ws.Cells.Item(6, 10).Value = 220  # was 219
$ws.Cells.Item(7, 10).Value = 396  # was 394
$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(2, 10).Value = 29  # was 27
$ws.Cells.Item(7, 10).Value = 83  # was 81
$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(3, 10).Value = 21  # was 20
$ws.Cells.Item(6, 10).Value = 36  # was 35
$ws.Cells.Item(7, 10).Value = 88  # was 86
$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(3, 10).Value = 82  # was 81
$ws.Cells.Item(6, 10).Value = 64  # was 63
$ws.Cells.Item(7, 10).Value = 240  # was 238
$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 10).Value = 246  # was 245
$ws.Cells.Item(6, 10).Value = 217  # was 216
$ws.Cells.Item(7, 10).Value = 730  # was 728
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(2, 10).Value = 48  # was 47
$ws.Cells.Item(6, 10).Value = 60  # was 59
$ws.Cells.Item(7, 10).Value = 172  # was 170
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(3, 10).Value = 113  # was 112
$ws.Cells.Item(6, 10).Value = 107  # was 106
$ws.Cells.Item(7, 10).Value = 353  # was 351
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(2, 10).Value = 28  # was 27
$ws.Cells.Item(7, 10).Value = 107  # was 106
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(6, 10).Value = 90  # was 89
$ws.Cells.Item(7, 10).Value = 192  # was 191
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 10).Value = 88  # was 86
$ws.Cells.Item(7, 10).Value = 308  # was 306
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(6, 10).Value = 123  # was 119
$ws.Cells.Item(7, 10).Value = 192  # was 188
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 10).Value = 82  # was 81
$ws.Cells.Item(6, 10).Value = 89  # was 88
$ws.Cells.Item(7, 10).Value = 250  # was 248
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(6, 10).Value = 76  # was 75
$ws.Cells.Item(7, 10).Value = 204  # was 203
$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(6, 10).Value = 164  # was 163
$ws.Cells.Item(7, 10).Value = 238  # was 237
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(6, 10).Value = 57  # was 56
$ws.Cells.Item(7, 10).Value = 153  # was 152
$ws = $wb.Worksheets.Item('Pullman')
$ws.Cells.Item(2, 10).Value = 34  # was 33
$ws.Cells.Item(7, 10).Value = 79  # was 78
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 10).Value = 97  # was 96
$ws.Cells.Item(3, 10).Value = 76  # was 75
$ws.Cells.Item(7, 10).Value = 276  # was 274
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(3, 10).Value = 85  # was 83
$ws.Cells.Item(6, 10).Value = 126  # was 125
$ws.Cells.Item(7, 10).Value = 316  # was 313
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(6, 10).Value = 45  # was 44
$ws.Cells.Item(7, 10).Value = 150  # was 149
$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(3, 10).Value = 30  # was 29
$ws.Cells.Item(7, 10).Value = 100  # was 99
$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Cells.Item(6, 10).Value = 11  # was 10
$ws.Cells.Item(7, 10).Value = 17  # was 16
